# Daily attendance processing - 2026-01-26 01:48:02
# Swap the order of the "Recorded By" contributors in column G:
#   "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# Only cells containing that exact combined value are touched; rows that
# were recorded solely by "dnasr281@gmail.com" or solely by "System" are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = $newValue
    }
}
